$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New FedEx job/tracking numbers that replace the previous 10 entries shown
# in column B (rows 2-11). These were appended to the shared-string table
# by the 1st-Apr-2022 test run; only the most recent 10 are surfaced on the
# visible report grid.
$newValues = @(
    "32323069",
    "32323070",
    "32323071",
    "32323072",
    "32323073",
    "32323074",
    "32323075",
    "32323076",
    "32323077",
    "32323078"
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 2)
    # Leading apostrophe forces text storage (these are id-like strings, not
    # numeric values) without leaving the cell's number format changed --
    # reset the style right after so no stray formatting sticks around.
    $cell.Value = "'" + $newValues[$i]
    $cell.Style = "Normal"
}
